$d = $word.ActiveDocument

# --- 1) Mark every inline picture's run as "no proofing" (<w:rPr><w:noProof/></w:rPr>) ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = $true
}

# --- 2) Extend the "Issue with TimeZone" bullet with " in Weather Underground" ---
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Issue with TimeZone") {
        $target = $d.Range($para.Range.Start, $para.Range.End - 1)
        $target.InsertAfter(" in Weather Underground")
        $found = $i
        break
    }
}

# --- 3) Insert a brand-new, completely empty paragraph right after that bullet ---
if ($found) {
    $nextPara = $d.Paragraphs.Item($found + 1)
    $insertPoint = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)
    [void]$insertPoint.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
}
